$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.695.19'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.059.46'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.75'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.31'
$ws.Range('E8').Value = '  -5.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.55'
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('E11').Value = '  -2.69%  '
$ws.Range('E12').Value = '  -3.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.939'
$ws.Range('E13').Value = '  +5.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.80'
$ws.Range('E14').Value = '  -3.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.360.40'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('E16').Value = '  -4.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.067.64'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.626.08'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  -6.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0868'
$ws.Range('E21').Value = '  -2.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.45'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.38'
$ws.Range('E25').Value = '  -2.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.32'
$ws.Range('E27').Value = '  -3.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.06'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.12'
$ws.Range('E31').Value = '  -8.19%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.20'
$ws.Range('E32').Value = '  +6.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.53'
$ws.Range('E33').Value = '  -4.64%  '
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0849'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('E39').Value = '  -4.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.00'
$ws.Range('E40').Value = '  -4.73%  '
$ws.Range('E41').Value = '  -6.39%  '
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('E43').Value = '  -3.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '94.84'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.413.71'
$ws.Range('E45').Value = '  +8.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0906'
$ws.Range('E46').Value = '  -8.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.70'
$ws.Range('E47').Value = '  +13.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.12'
$ws.Range('E48').Value = '  -5.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.92'
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('E50').Value = '  -4.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.249.42'
$ws.Range('E51').Value = '  +0.93%  '
